$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oktober")

$ws.Range("B26").Value = 1
$ws.Range("B27").Value = 2

$ws.Range("C27").Value = "DMX und ESP ferigstellen"
$ws.Range("C26").Value = "DMX und ESP weitere arbeit"
$ws.Range("D26").Value = "11:00 - 12:00"
$ws.Range("D27").Value = "14:00 - 16:00"

$ws.Activate()
$ws.Range("D28").Select()
